$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Good Morning" with "GIT UPDATE" in cell E8 (this also removes the
# now-unused "Good Morning" shared string and appends "GIT UPDATE" as a new
# shared string entry once the workbook is saved).
$ws.Range("E8").Value = "GIT UPDATE"

# Match the saved selection state (sheet view shows E8 as the active cell).
$ws.Range("E8").Select()
